# Penambahan input kecamatan, kelurahan, kota dan kode pos di module Customer
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customer")

# Current header row (A1:E1) = NAMA | ALAMAT | KONTAK | TELEPON | PLAFON PIUTANG
# We need to insert 4 new columns (KECAMATAN, KELURAHAN, KOTA, KODE POS) after
# ALAMAT (column B) and before KONTAK (column C), shifting the old C:E to G:I.

# Insert four new blank columns before the current column C (KONTAK).
$ws.Range("C1:F1").EntireColumn.Insert()

# Fill in the new header cells (order matches the shared-string insertion
# order recorded in the original commit: KELURAHAN, KECAMATAN, KOTA, KODE POS).
$ws.Range("D1").Value = "KELURAHAN"
$ws.Range("C1").Value = "KECAMATAN"
$ws.Range("E1").Value = "KOTA"
$ws.Range("F1").Value = "KODE POS"

# Copy the header style used by the other header cells (e.g. B1) onto the
# newly inserted header cells so they match the existing formatting.
$ws.Range("B1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new column widths for the inserted columns and the shifted ones.
$ws.Range("C1:E1").ColumnWidth = 16.28515625
$ws.Range("F1").ColumnWidth = 9.85546875

# Update the active selection to match the target state.
$ws.Range("B8").Select()
